$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About": rework the Notes section (rows 17-28)
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Row 17 keeps the "Notes:" header text (same text, new shared-string slot).
$about.Range("A17").Value = "Notes:"

# Row 23-25: the old "motorbikes" note text is replaced by the LCFS note.
# (Set first so the new shared strings land in the same order as the source file.)
$about.Range("A23").Value = "For the LCFS to work correctly, all on-road modes need to have the same value, so we"
$about.Range("A24").Value = "use the passenger LDVs value for all on-road modes.  (It is very close to the calculated"
$about.Range("A25").Value = "passenger HDVs value in any case.)"

# Rows 18-21: new explanatory note, using the same (slightly distinct) applied
# font style as the rest of the note block.
$about.Range("A18").Value = "For vehicle types that can use electricity, this variable specifies the percentage"
$about.Range("A19").Value = "reduction in fuel use (on a BTU basis) relative to the typical fuel type for that vehicle"
$about.Range("A20").Value = "type (e.g. gasoline for LDVs, diesel for HDVs, etc.) due to the fact that electricity"
$about.Range("A21").Value = "can be converted into work more efficiently than other fuel types."

$about.Range("A18:A22").Font.ThemeFont = $about.Range("A1").Font.ThemeFont

# ---------------------------------------------------------------------------
# Sheet "PTFURfE": widen column A, retitle header, rework the B/C formulas
# ---------------------------------------------------------------------------
$ptf = $wb.Worksheets.Item("PTFURfE")

$ptf.Columns.Item(1).ColumnWidth = 16

$ptf.Range("A1").Value = "Percentage Reduction (dimensionless)"
$ptf.Range("A1").Font.Bold = $true
$ptf.Range("A1").WrapText = $true
$ptf.Rows.Item(1).RowHeight = 45

# Row 27-28 (About sheet): new note about aircraft/ships -- set after the
# PTFURfE header above so the new shared strings land in source-file order.
$about.Range("A27").Value = "Aircraft and ships are assumed to be the same as rail, since they all use large engines"
$about.Range("A28").Value = "intended to move heavy craft."

# Row 3 (HDVs): both Passengers & Freight now simply mirror the LDVs value.
$ptf.Range("B3").Formula = '=$B$2'
$ptf.Range("C3").NumberFormat = $ptf.Range("B2").NumberFormat
$ptf.Range("C3").Formula = '=$B$2'

# Row 4 (aircraft): mirrors row 5 (rail) instead of being hard-coded zero.
$ptf.Range("B4").NumberFormat = $ptf.Range("B2").NumberFormat
$ptf.Range("B4").Formula = '=B5'
$ptf.Range("C4").NumberFormat = $ptf.Range("B2").NumberFormat
$ptf.Range("C4").Formula = '=C5'

# Row 5 (rail): Freight now mirrors Passengers.
$ptf.Range("C5").NumberFormat = $ptf.Range("B2").NumberFormat
$ptf.Range("C5").Formula = '=B5'

# Row 6 (ships): mirrors row 5 (rail) instead of being hard-coded zero.
$ptf.Range("B6").NumberFormat = $ptf.Range("B2").NumberFormat
$ptf.Range("B6").Formula = '=B5'
$ptf.Range("C6").NumberFormat = $ptf.Range("B2").NumberFormat
$ptf.Range("C6").Formula = '=C5'

# Row 7 (motorbikes): anchor the LDVs reference absolutely and mirror it to Freight.
$ptf.Range("B7").Formula = '=$B$2'
$ptf.Range("C7").NumberFormat = $ptf.Range("B2").NumberFormat
$ptf.Range("C7").Formula = '=$B$2'
